# HMI Saw RUN Updated
#
# 1) Touch up the view state on the two existing sheets that lost focus
#    (NotesForSeptVisit, HMI Screen Shots) so their selections match the
#    committed state. Doing this before the new sheet is created/activated
#    means they naturally lose "tabSelected" once the new sheet becomes
#    the active tab.
$wb = $excel.ActiveWorkbook

$wsNotes = $wb.Worksheets.Item("NotesForSeptVisit")
$wsNotes.Range("C3").Select()

$wsShots = $wb.Worksheets.Item("HMI Screen Shots")
$wsShots.Range("N16").Select()

# 2) Add the new "Sheet1" tab after the last existing sheet (HMI Screen
#    Shots) so it becomes sheet #6 / sheetId 6 and the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# 3) Populate the RUN/OLD-NUMBER-NEW mapping table.
#    Values are written in the exact order the strings were first used
#    so the shared-string table comes out in the same sequence as the
#    authored workbook.

$ws.Range("C3").Value = "Control Power On"
$ws.Range("C5").Value = "Loop Pump Aux"
$ws.Range("C4").Value = "Pump Motor Aux"
$ws.Range("C6").Value = "Cooling Fan Aux"
$ws.Range("C8").Value = "Swarf Pump Aux"
$ws.Range("C9").Value = "Swarf Conveyor Aux"
$ws.Range("C19").Value = "Gauge Home Limit"
$ws.Range("C20").Value = "Traverse Home Limit"
$ws.Range("C21").Value = "Raise Lower Home Limit"
$ws.Range("C22").Value = "Guage In OverTravel"
$ws.Range("C23").Value = "Gauge Out Overtravel"
$ws.Range("C24").Value = "Traverse Forward Overtravel"
$ws.Range("C25").Value = "Traverse Reverse Overtravel"
$ws.Range("C26").Value = "RaiseOvertravel"
$ws.Range("C27").Value = "LowerOvertravel"
$ws.Range("G10").Value = "HPU Motor Starter"
$ws.Range("C35").Value = "HPU Motor Starter"
$ws.Range("G11").Value = "Loop Motor Starter"
$ws.Range("C36").Value = "Loop Motor Starter"
$ws.Range("G14").Value = "Swarf Motor Starter"
$ws.Range("C39").Value = "Swarf Motor Starter"
$ws.Range("G15").Value = "Swarf Conveyor Starter"
$ws.Range("C40").Value = "Swarf Conveyor Starter"
$ws.Range("G20").Value = "Laser Line CR"
$ws.Range("C45").Value = "Laser Line CR"
$ws.Range("G21").Value = "Lt Auto On Green"
$ws.Range("C46").Value = "Lt Auto On Green"
$ws.Range("G22").Value = "Lt Power On Amber"
$ws.Range("C47").Value = "Lt Power On Amber"
$ws.Range("G23").Value = "Lt Red Axis Active"
$ws.Range("C48").Value = "Lt Red Axis Active"
$ws.Range("B1").Value = "OLD"
$ws.Range("A2").Value = "NUMBER"
$ws.Range("B2").Value = "NUMBER"
$ws.Range("A1").Value = "NEW"
$ws.Range("C14").Value = "Blade Guard"
$ws.Range("C15").Value = "Blade Detect"
$ws.Range("C16").Value = "Swarf Filter"
$ws.Range("G12").Value = "Hyd Cooling"
$ws.Range("C37").Value = "Hyd Cooling"
$ws.Range("G13").Value = "Coolant Starter"
$ws.Range("C38").Value = "Coolant Starter"
# 4) Fill in the numeric OLD/NEW tag-number columns (A, B, E, F).
$ws.Range("A3").Value = 502
$ws.Range("B3").Value = 600
$ws.Range("F3").Value = 625
$ws.Range("A4").Value = 503
$ws.Range("B4").Value = 601
$ws.Range("F4").Value = 626
$ws.Range("A5").Value = 504
$ws.Range("B5").Value = 602
$ws.Range("F5").Value = 627
$ws.Range("A6").Value = 505
$ws.Range("B6").Value = 603
$ws.Range("F6").Value = 628
$ws.Range("B7").Value = 604
$ws.Range("F7").Value = 629
$ws.Range("A8").Value = 534
$ws.Range("B8").Value = 605
$ws.Range("F8").Value = 630
$ws.Range("A9").Value = 535
$ws.Range("B9").Value = 606
$ws.Range("F9").Value = 631
$ws.Range("B10").Value = 607
$ws.Range("E10").Value = 614
$ws.Range("F10").Value = 632
$ws.Range("B11").Value = 608
$ws.Range("E11").Value = 615
$ws.Range("F11").Value = 633
$ws.Range("B12").Value = 609
$ws.Range("E12").Value = 616
$ws.Range("F12").Value = 634
$ws.Range("B13").Value = 610
$ws.Range("E13").Value = 617
$ws.Range("F13").Value = 635
$ws.Range("A14").Value = 545
$ws.Range("B14").Value = 611
$ws.Range("E14").Value = 630
$ws.Range("F14").Value = 636
$ws.Range("A15").Value = 602
$ws.Range("B15").Value = 612
$ws.Range("E15").Value = 631
$ws.Range("F15").Value = 637
$ws.Range("A16").Value = 603
$ws.Range("B16").Value = 613
$ws.Range("F16").Value = 638
$ws.Range("B17").Value = 614
$ws.Range("F17").Value = 639
$ws.Range("B18").Value = 615
$ws.Range("F18").Value = 640
$ws.Range("A19").Value = 516
$ws.Range("B19").Value = 616
$ws.Range("F19").Value = 641
$ws.Range("A20").Value = 517
$ws.Range("B20").Value = 617
$ws.Range("E20").Value = 642
$ws.Range("F20").Value = 642
$ws.Range("A21").Value = 530
$ws.Range("B21").Value = 618
$ws.Range("E21").Value = 643
$ws.Range("F21").Value = 643
$ws.Range("A22").Value = 531
$ws.Range("B22").Value = 619
$ws.Range("E22").Value = 644
$ws.Range("F22").Value = 644
$ws.Range("A23").Value = 532
$ws.Range("B23").Value = 620
$ws.Range("E23").Value = 645
$ws.Range("F23").Value = 645
$ws.Range("A24").Value = 533
$ws.Range("B24").Value = 621
$ws.Range("F24").Value = 646
$ws.Range("A25").Value = 542
$ws.Range("B25").Value = 622
$ws.Range("A26").Value = 543
$ws.Range("B26").Value = 623
$ws.Range("A27").Value = 544
$ws.Range("B27").Value = 624
$ws.Range("B28").Value = 625
$ws.Range("B29").Value = 626
$ws.Range("B30").Value = 627
$ws.Range("B31").Value = 628
$ws.Range("B32").Value = 629
$ws.Range("B33").Value = 630
$ws.Range("B34").Value = 631
$ws.Range("A35").Value = 614
$ws.Range("B35").Value = 632
$ws.Range("A36").Value = 615
$ws.Range("B36").Value = 633
$ws.Range("A37").Value = 616
$ws.Range("B37").Value = 634
$ws.Range("A38").Value = 617
$ws.Range("B38").Value = 635
$ws.Range("A39").Value = 630
$ws.Range("B39").Value = 636
$ws.Range("A40").Value = 631
$ws.Range("B40").Value = 637
$ws.Range("B41").Value = 638
$ws.Range("B42").Value = 639
$ws.Range("B43").Value = 640
$ws.Range("B44").Value = 641
$ws.Range("A45").Value = 642
$ws.Range("B45").Value = 642
$ws.Range("A46").Value = 643
$ws.Range("B46").Value = 643
$ws.Range("A47").Value = 644
$ws.Range("B47").Value = 644
$ws.Range("A48").Value = 645
$ws.Range("B48").Value = 645
$ws.Range("B49").Value = 646

# 5) Column C ("OLD"/new tag name) is sized to fit its longest label.
$ws.Columns.Item(3).ColumnWidth = 15.6

# 6) Leave the selection on this (now active) sheet matching the
#    authored workbook: E3:G24 with the active cell on E3.
$ws.Range("E3:G24").Select()
